# Converts an EMU (English Metric Unit) value to points for use with the
# Shape.Left / Shape.Top / Shape.Width / Shape.Height properties (which are
# Single-precision floats under the hood). A tiny epsilon is added so that
# the float32 round-trip back to EMU lands on the intended integer value
# instead of being truncated one EMU short.
function EmuToPt($emu) {
    return ($emu / 914400.0 * 72.0) + 0.00001
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "27" (the "2" step-number tile) : move down -------------------
$shpStepNumber = $s.Shapes.Item(8)
$shpStepNumber.Left = EmuToPt(218936)
$shpStepNumber.Top = EmuToPt(3243672)

# --- Shape "28" ("Criteria for success" heading) : reposition ------------
$shpCriteria = $s.Shapes.Item(9)
$shpCriteria.Left = EmuToPt(588238)
$shpCriteria.Top = EmuToPt(3295840)

# --- Shape "37" (Challenges bullet list) : bump font size 9.5 -> 9.75 pt -
$shpChallenges = $s.Shapes.Item(18)
$trChallenges = $shpChallenges.TextFrame.TextRange
for ($i = 1; $i -le $trChallenges.Paragraphs().Count; $i++) {
    $trChallenges.Paragraphs($i).Font.Size = 9.75
}

# --- Shape "38" (Data source paragraph) : expand wording ------------------
$shpDataSource = $s.Shapes.Item(19)
$shpDataSource.TextFrame.TextRange.Text = "The data source for this project is a CSV file provided by the Database Manager, Alesha Eisen. This file includes data on 330 resorts in the same market segment as Big Mountain Resort. This data, along with the metadata file with column descriptions, is used to identify important variables for our analysis."

# --- Shape "47" (Key stakeholders list) : bump font size 9.5 -> 9.75 pt --
$shpStakeholders = $s.Shapes.Item(28)
$trStakeholders = $shpStakeholders.TextFrame.TextRange
for ($i = 1; $i -le $trStakeholders.Paragraphs().Count; $i++) {
    $trStakeholders.Paragraphs($i).Font.Size = 9.75
}

Write-Output "Edits applied"
